# Testing + model polishing
# Updated syngas production and upstream emissions
# New testing file with HIsarna energy mix.
# Disconnected coke from pellet production.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuels")

# --- Disconnect the coke-from-pellet-production formulas: D10 and D15
#     go from a live formula to a plain literal 0 ---
$ws.Range("D10").Value = 0
$ws.Range("D15").Value = 0

# --- New testing file with HIsarna energy mix: row 16 (syngas - wood)
#     HHV/LHV nudged from 21.6 to 21.5 ---
$ws.Range("B16").Value = 21.5
$ws.Range("C16").Value = 21.5

# --- Note on how the new HIsarna number was derived ---
$cmt = $ws.Range("A16").AddComment("Microsoft Office User:" + [char]10 + "calculated from Swanson 2010 (HT scenario)")

# --- Move the saved cursor / selection to where the author left off ---
$ws.Range("D17").Select() | Out-Null
